# cryptos.xlsx refresh -- GitHub Actions scheduled data pull.
# Rewrites the Coin / Link / Price / Volume(1h) table (rows 2-51) with the
# latest values. The coin list shifted by one row in the middle of the
# table (a new coin -- WrappedliquidstakedEther2.0 -- was inserted at row
# 14, pushing everything below it down by one and dropping the former last
# row, Algorand, off the bottom), so Coin/Link/Price/Volume are all
# rewritten for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores numeric-looking text (e.g. "26.057.46",
# "10.07", "7.220"). Force Text format before writing so Excel does not
# auto-convert it to a Number (which would also strip significant
# trailing zeros, e.g. "7.220" -> 7.22).
$priceCells = @(
    "D2", "D3", "D5", "D6", "D9", "D10", "D12", "D13", "D14", "D15",
    "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25",
    "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35",
    "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45",
    "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($c in $priceCells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range("D2").Value = "26.057.46"
$ws.Range("D3").Value = "1.651.04"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "217.22"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "0.5263"
$ws.Range("E6").Value = "  +1.92%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("D9").Value = "0.06323"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").Value = "20.34"
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "4.501"
$ws.Range("D13").Value = "1.620.80"
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.879.01"
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.5487"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0₅8205"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "65.46"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "26.070.12"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "4.573"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "190.54"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "10.07"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "6.021"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "143.58"
$ws.Range("E25").Value = "  +2.73%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "0.1237"
$ws.Range("E26").Value = "  +1.29%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "7.220"
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "16.04"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "1.427"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "0.05813"
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.274"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "3.546"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "3.269"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "1.580"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "0.9458"
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "2.781"
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.409"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "0.5730"
$ws.Range("E38").Value = "  +1.09%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01609"
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "0.8441"
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "5.742"
$ws.Range("E41").Value = "  -5.02%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "103.69"
$ws.Range("E43").Value = "  +3.00%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.029.05"
$ws.Range("E44").Value = "  +1.78%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.795.08"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "56.84"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.4323"
$ws.Range("E48").Value = "  +2.90%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.857"
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05144"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "1.460"
$ws.Range("E51").Value = "  +1.18%  "
